# Add a new column (BP, col 68) with header "13-sep" and per-product counts,
# mirroring the existing BO (col 67) "12-sep" column that sits right before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCol = 68
$sourceCol = 67

# Header cell (row 1) - new date label, copy the header's style/format
$headerCell = $ws.Cells.Item(1, $newCol)
$headerCell.Value = "13-sep"
$headerCell.NumberFormat = $ws.Cells.Item(1, $sourceCol).NumberFormat

# Data values for rows 2-11, copying number format/alignment from the
# neighboring BO column so the new cells match the existing style.
$values = @{
    2 = 14
    3 = 15
    4 = 12
    5 = 14
    6 = 10
    7 = 17
    8 = 21
    9 = 10
    10 = 10
    11 = 14
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, $newCol)
    $cell.Value = $values[$row]
    $cell.HorizontalAlignment = $ws.Cells.Item($row, $sourceCol).HorizontalAlignment
    $cell.NumberFormat = $ws.Cells.Item($row, $sourceCol).NumberFormat
}

# Update the active selection to match the saved workbook state.
$ws.Range("BR8").Select() | Out-Null
